$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.161.51"
Set-TextValue "E2" "  -0.09%  "
Set-TextValue "D3" "1.652.07"
Set-TextValue "E3" "  -0.42%  "
Set-TextValue "E4" "  +0.49%  "
Set-TextValue "D5" "217.99"
Set-TextValue "E5" "  +1.03%  "
Set-TextValue "D6" "0.5198"
Set-TextValue "E6" "  -1.24%  "
Set-TextValue "E7" "  +0.44%  "
Set-TextValue "D8" "0.2644"
Set-TextValue "E8" "  +0.46%  "
Set-TextValue "D9" "0.06310"
Set-TextValue "E9" "  -1.47%  "
Set-TextValue "D10" "21.19"
Set-TextValue "E10" "  +1.15%  "
Set-TextValue "D11" "0.07732"
Set-TextValue "E11" "  -0.49%  "
Set-TextValue "D12" "1.647.27"
Set-TextValue "E12" "  -0.68%  "
Set-TextValue "D13" "4.419"
Set-TextValue "E13" "  -1.06%  "
Set-TextValue "D14" "0.5453"
Set-TextValue "E14" "  -1.61%  "
Set-TextValue "D15" "0.0₅8198"
Set-TextValue "E15" "  -1.13%  "
Set-TextValue "D16" "64.67"
Set-TextValue "E16" "  -1.00%  "
Set-TextValue "D17" "26.183.81"
Set-TextValue "E17" "  -0.03%  "
Set-TextValue "E18" "  +0.41%  "
Set-TextValue "E19" "  -1.84%  "
Set-TextValue "D20" "191.54"
Set-TextValue "E20" "  +0.37%  "
Set-TextValue "D21" "10.15"
Set-TextValue "E21" "  -1.54%  "
Set-TextValue "D22" "6.166"
Set-TextValue "E22" "  -3.25%  "
Set-TextValue "D23" "1.008"
Set-TextValue "E23" "  +0.64%  "
Set-TextValue "D24" "138.48"
Set-TextValue "E24" "  -3.16%  "
Set-TextValue "D25" "0.1239"
Set-TextValue "E25" "  -1.76%  "
Set-TextValue "D26" "7.277"
Set-TextValue "E26" "  -2.02%  "
Set-TextValue "D27" "16.04"
Set-TextValue "E27" "  -0.02%  "
Set-TextValue "E28" "  -1.05%  "
Set-TextValue "D29" "0.06056"
Set-TextValue "E29" "  -1.74%  "
Set-TextValue "D30" "1.282"
Set-TextValue "E30" "  +1.22%  "
Set-TextValue "E31" "  -0.53%  "
Set-TextValue "D32" "3.356"
Set-TextValue "E32" "  -2.14%  "
Set-TextValue "D33" "1.654"
Set-TextValue "D34" "0.9830"
Set-TextValue "E34" "  -1.86%  "
Set-TextValue "D35" "2.412"
Set-TextValue "E35" "  +0.47%  "
Set-TextValue "E36" "  +0.26%  "
Set-TextValue "D37" "0.5936"
Set-TextValue "E37" "  +4.43%  "
Set-TextValue "E38" "  -0.68%  "
Set-TextValue "D39" "5.958"
Set-TextValue "E39" "  +0.66%  "
Set-TextValue "D40" "0.8633"
Set-TextValue "E40" "  +0.97%  "
Set-TextValue "D41" "1.052.75"
Set-TextValue "E41" "  +1.87%  "
Set-TextValue "D42" "1.004"
Set-TextValue "E42" "  +0.23%  "
Set-TextValue "D43" "99.69"
Set-TextValue "E43" "  +0.05%  "
Set-TextValue "D44" "1.794.26"
Set-TextValue "E44" "  -0.69%  "
Set-TextValue "B45" "BabyDogeCoin"
Set-TextValue "C45" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D45" "0.0₈112"
Set-TextValue "E45" "  +4.61%  "
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "57.34"
Set-TextValue "E46" "  +1.99%  "
Set-TextValue "D47" "1.005"
Set-TextValue "E47" "  +0.16%  "
Set-TextValue "D48" "8.055"
Set-TextValue "E48" "  -0.18%  "
Set-TextValue "D49" "0.05175"
Set-TextValue "E49" "  +0.14%  "
Set-TextValue "B50" "Mantle"
Set-TextValue "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D50" "0.4233"
Set-TextValue "E50" "  +0.47%  "
Set-TextValue "B51" "RenderToken"
Set-TextValue "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "1.462"
Set-TextValue "E51" "  +4.50%  "
